$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# 1) "_settings" sheet: rename the "dataType" parameter row to "dateType"
#    (both the parameter key in column A and the display name in column B,
#    row 12 of the paramTable).
# ---------------------------------------------------------------------------
$settings = $wb.Worksheets.Item("_settings")
$settings.Range("A12").Value = "dateType"
$settings.Range("B12").Value = "dateType"

# ---------------------------------------------------------------------------
# 2) "_input" sheet: guard the D11 formula so a negative delta clamps to 0.
# ---------------------------------------------------------------------------
$input = $wb.Worksheets.Item("_input")
$input.Range("D11").Formula = "=IF(`$E3<0,0,`$E3)"

# ---------------------------------------------------------------------------
# 3) "trend" sheet: two new blank marker cells (I1 / I23) and widen the used
#    range, then restore row 1's explicit 3pt height (writing into I1 makes
#    Excel recompute an autofit height otherwise).
# ---------------------------------------------------------------------------
$trend = $wb.Worksheets.Item("trend")
$trend.Range("I1").Value = "  "
$trend.Range("I23").Value = "  "
$trend.Rows.Item(1).RowHeight = 3

# ---------------------------------------------------------------------------
# 4) "trend" sheet: reposition/resize the embedded "trendChart" chart.
#    New anchor: from col0/row3 (colOff 38098, rowOff 47626) to col7/row21
#    (colOff 390525, rowOff 133350). Convert to points (Left/Top/Width/
#    Height) using the sheet's current column widths / row heights.
# ---------------------------------------------------------------------------
$chartObj = $trend.ChartObjects().Item(1)

$emuPerPoint = 12700

# from-anchor: column 0 (i.e. 0 whole columns before it), row index 3 (3 whole rows above it)
$leftPt = 0 + (38098 / $emuPerPoint)

$topPt = 0
for ($r = 1; $r -le 3; $r++) { $topPt += $trend.Cells.Item($r, 1).Height }
$topPt += 47626 / $emuPerPoint

# to-anchor: 7 whole columns before it (A..G), 21 whole rows above it
$rightPt = 0
for ($c = 1; $c -le 7; $c++) { $rightPt += $trend.Cells.Item(1, $c).Width }
$rightPt += 390525 / $emuPerPoint

$bottomPt = 0
for ($r = 1; $r -le 21; $r++) { $bottomPt += $trend.Cells.Item($r, 1).Height }
$bottomPt += 133350 / $emuPerPoint

$chartObj.Left = $leftPt
$chartObj.Top = $topPt
$chartObj.Width = $rightPt - $leftPt
$chartObj.Height = $bottomPt - $topPt
